$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45973.4143287037
$ws.Range("D3").Value = 45973.4143287037
$ws.Range("D4").Value = 45973.4143287037
$ws.Range("D5").Value = 45973.4143287037
$ws.Range("D6").Value = 45973.4143287037
$ws.Range("D7").Value = 45973.4143287037
$ws.Range("D8").Value = 45973.4143287037
$ws.Range("D9").Value = 45973.4143287037
$ws.Range("D10").Value = 45973.4143287037
$ws.Range("D11").Value = 45973.4143287037
$ws.Range("D12").Value = 45973.4143287037
$ws.Range("D13").Value = 45973.4143287037
$ws.Range("D14").Value = 45973.4143287037
$ws.Range("D15").Value = 45973.4143287037
$ws.Range("D16").Value = 45973.4143287037
$ws.Range("D17").Value = 45973.4143287037
$ws.Range("D18").Value = 45973.4143287037

$ws.Range("A19").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B19").Value = "602号直流"
$ws.Range("C19").Value = 45971.191238425927
$ws.Range("D19").Value = 45973.4143287037
$ws.Range("A20").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B20").Value = "705号直流"
$ws.Range("C20").Value = 45971.419039351851
$ws.Range("D20").Value = 45973.4143287037
$ws.Range("A21").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B21").Value = "406号直流"
$ws.Range("C21").Value = 45971.54614583333
$ws.Range("D21").Value = 45973.4143287037
$ws.Range("A22").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B22").Value = "503号直流"
$ws.Range("C22").Value = 45971.561631944445
$ws.Range("D22").Value = 45973.4143287037
$ws.Range("A23").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B23").Value = "101号直流"
$ws.Range("C23").Value = 45971.970555555556
$ws.Range("D23").Value = 45973.4143287037
$ws.Range("A24").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B24").Value = "505号直流"
$ws.Range("C24").Value = 45972.035127314812
$ws.Range("D24").Value = 45973.4143287037
$ws.Range("A25").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B25").Value = "A01号直流"
$ws.Range("C25").Value = 45972.038611111115
$ws.Range("D25").Value = 45973.4143287037
$ws.Range("A26").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B26").Value = "703号直流"
$ws.Range("C26").Value = 45972.050416666665
$ws.Range("D26").Value = 45973.4143287037
$ws.Range("A27").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B27").Value = "702号直流"
$ws.Range("C27").Value = 45972.123761574076
$ws.Range("D27").Value = 45973.4143287037
$ws.Range("A28").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B28").Value = "401号直流"
$ws.Range("C28").Value = 45972.159884259258
$ws.Range("D28").Value = 45973.4143287037
$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "904号直流"
$ws.Range("C29").Value = 45972.241493055553
$ws.Range("D29").Value = 45973.4143287037
$ws.Range("A30").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B30").Value = "208号直流"
$ws.Range("C30").Value = 45972.507245370369
$ws.Range("D30").Value = 45973.4143287037
$ws.Range("A31").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B31").Value = "B03号直流"
$ws.Range("C31").Value = 45972.508888888886
$ws.Range("D31").Value = 45973.4143287037
$ws.Range("A32").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B32").Value = "306号直流"
$ws.Range("C32").Value = 45972.537523148145
$ws.Range("D32").Value = 45973.4143287037
$ws.Range("A33").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B33").Value = "107号直流"
$ws.Range("C33").Value = 45972.544791666667
$ws.Range("D33").Value = 45973.4143287037
$ws.Range("A34").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B34").Value = "504号直流"
$ws.Range("C34").Value = 45972.562094907407
$ws.Range("D34").Value = 45973.4143287037
$ws.Range("A35").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B35").Value = "203号直流"
$ws.Range("C35").Value = 45972.569097222222
$ws.Range("D35").Value = 45973.4143287037
$ws.Range("A36").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B36").Value = "103号直流"
$ws.Range("C36").Value = 45972.572881944441
$ws.Range("D36").Value = 45973.4143287037
$ws.Range("A37").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B37").Value = "403号直流"
$ws.Range("C37").Value = 45972.574583333335
$ws.Range("D37").Value = 45973.4143287037
$ws.Range("A38").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B38").Value = "105号直流"
$ws.Range("C38").Value = 45972.575798611113
$ws.Range("D38").Value = 45973.4143287037
$ws.Range("A39").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B39").Value = "002B号直流"
$ws.Range("C39").Value = 45972.585972222223
$ws.Range("D39").Value = 45973.4143287037
$ws.Range("A40").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B40").Value = "402号直流"
$ws.Range("C40").Value = 45972.586342592593
$ws.Range("D40").Value = 45973.4143287037
$ws.Range("A41").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B41").Value = "802号直流"
$ws.Range("C41").Value = 45972.588275462964
$ws.Range("D41").Value = 45973.4143287037
$ws.Range("A42").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B42").Value = "501号直流"
$ws.Range("C42").Value = 45972.600011574075
$ws.Range("D42").Value = 45973.4143287037
$ws.Range("A43").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B43").Value = "402号直流"
$ws.Range("C43").Value = 45972.631921296299
$ws.Range("D43").Value = 45973.4143287037
$ws.Range("A44").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B44").Value = "008A号直流"
$ws.Range("C44").Value = 45972.644895833335
$ws.Range("D44").Value = 45973.4143287037
$ws.Range("A45").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B45").Value = "A04号直流"
$ws.Range("C45").Value = 45972.654374999998
$ws.Range("D45").Value = 45973.4143287037
$ws.Range("A46").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B46").Value = "103号直流"
$ws.Range("C46").Value = 45972.656111111108
$ws.Range("D46").Value = 45973.4143287037
$ws.Range("A47").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B47").Value = "210号直流"
$ws.Range("C47").Value = 45972.659328703703
$ws.Range("D47").Value = 45973.4143287037
$ws.Range("A48").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B48").Value = "206号直流"
$ws.Range("C48").Value = 45972.684016203704
$ws.Range("D48").Value = 45973.4143287037
$ws.Range("A49").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B49").Value = "502号直流"
$ws.Range("C49").Value = 45972.685428240744
$ws.Range("D49").Value = 45973.4143287037
$ws.Range("A50").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B50").Value = "902号直流"
$ws.Range("C50").Value = 45972.712581018517
$ws.Range("D50").Value = 45973.4143287037
$ws.Range("A51").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B51").Value = "206号直流"
$ws.Range("C51").Value = 45972.724432870367
$ws.Range("D51").Value = 45973.4143287037
$ws.Range("A52").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B52").Value = "905号直流"
$ws.Range("C52").Value = 45972.750127314815
$ws.Range("D52").Value = 45973.4143287037
$ws.Range("A53").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B53").Value = "102号直流"
$ws.Range("C53").Value = 45972.758530092593
$ws.Range("D53").Value = 45973.4143287037
$ws.Range("A54").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B54").Value = "704号直流"
$ws.Range("C54").Value = 45972.855810185189
$ws.Range("D54").Value = 45973.4143287037

[void]$ws.Range("F12").Select()